$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# Update date values for rows 8-11 (columns B, C, Q, R)
# B: inicio periodo, C: fin periodo -> shifted from Q2 2022 to Q3 2022
# Q, R: fecha de validacion -> shifted forward as well
foreach ($r in 8..11) {
    $ws.Range("B$r").Value2 = 44743
    $ws.Range("C$r").Value2 = 44834
    $ws.Range("Q$r").Value2 = 44844
    $ws.Range("R$r").Value2 = 44844
}

# Update row heights for rows 8-11 from 90 to 60
foreach ($r in 8..11) {
    $ws.Rows.Item($r).RowHeight = 60
}

# Update the active sheet view: scroll position and selection
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 6
$ws.Range("I14").Select()
